$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36 - this shifts all existing rows 36..154
# down to 37..155 (the row that falls off the used range, old row 154,
# becomes the new row 155), matching the target diff exactly.
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly data point.
$ws.Range("A36").Value = 11
$ws.Range("B36").Value = "Vega Monumental Concepción"
$ws.Range("C36").Value = "Bíobío"
$ws.Range("D36").Value = 45274
$ws.Range("E36").Value = 8
$ws.Range("F36").Value = 100112012
$ws.Range("G36").Value = "Espinaca"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 100
$ws.Range("K36").Value = 8000
$ws.Range("L36").Value = 8000
$ws.Range("M36").Value = 8000
$ws.Range("N36").Value = "$/cuna 10 kilos"
$ws.Range("O36").Value = "Región Metropolitana"
$ws.Range("P36").Value = 800
$ws.Range("Q36").Value = 10
$ws.Range("R36").Value = "Hortaliza"
